$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.448.15'
$ws.Range("E2").Value = '  -6.49%  '
$ws.Range("D3").Value = '2.198.02'
$ws.Range("E3").Value = '  -7.13%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.86'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.23'
$ws.Range("E6").Value = '  -11.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").Value = '  -9.41%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.550'
$ws.Range("E9").Value = '  -10.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.10'
$ws.Range("E10").Value = '  -11.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.55'
$ws.Range("E11").Value = '  -3.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0819'
$ws.Range("E12").Value = '  -10.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.51'
$ws.Range("E13").Value = '  -11.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.104'
$ws.Range("E14").Value = '  -4.41%  '
$ws.Range("D15").Value = '2.537.68'
$ws.Range("E15").Value = '  -7.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.843'
$ws.Range("E16").Value = '  -13.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.80'
$ws.Range("E17").Value = '  -10.84%  '
$ws.Range("D18").Value = '2.193.28'
$ws.Range("E18").Value = '  -8.01%  '
$ws.Range("D19").Value = '42.330.23'
$ws.Range("E19").Value = '  -6.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.05'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.43'
$ws.Range("E21").Value = '  -12.10%  '
$ws.Range("D22").Value = '0.0₃0938'
$ws.Range("E22").Value = '  -11.80%  '
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.17'
$ws.Range("E23").Value = '  -8.80%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.60'
$ws.Range("E24").Value = '  -12.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '231.86'
$ws.Range("E25").Value = '  -11.35%  '
$ws.Range("E26").Value = '  -9.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("E28").Value = '  -9.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.15'
$ws.Range("E29").Value = '  -8.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.38'
$ws.Range("E30").Value = '  -13.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.25'
$ws.Range("E31").Value = '  -9.89%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0862'
$ws.Range("E32").Value = '  -10.07%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '156.20'
$ws.Range("E33").Value = '  -8.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.07'
$ws.Range("E34").Value = '  -12.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.70'
$ws.Range("E35").Value = '  -8.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.14'
$ws.Range("E36").Value = '  +5.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.120'
$ws.Range("E37").Value = '  -7.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.37'
$ws.Range("E38").Value = '  -9.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.81'
$ws.Range("E39").Value = '  +4.95%  '
$ws.Range("E40").Value = '  -12.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.47'
$ws.Range("E41").Value = '  -12.90%  '
$ws.Range("E42").Value = '  -11.42%  '
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = '1.768.72'
$ws.Range("E44").Value = '  +7.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '88.51'
$ws.Range("E45").Value = '  -13.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.74'
$ws.Range("E46").Value = '  -11.20%  '
$ws.Range("E47").Value = '  -13.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '75.90'
$ws.Range("E48").Value = '  -6.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.25'
$ws.Range("E49").Value = '  -5.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '59.49'
$ws.Range("E50").Value = '  -14.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.39'
$ws.Range("E51").Value = '  -10.13%  '
